$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Replace the heading text "1.3.2. Sustav za upravljanje klijentima (CRM)"
#    with "1.3.2. Sustav za upravljanje klijentima " (drop the "(CRM)" marker,
#    keep the trailing space) and remember where the replaced run ends so we
#    can re-plant the _GoBack bookmark right after it.
# ---------------------------------------------------------------------------
$headingRng = $d.Content
$headingRng.Find.ClearFormatting()
$found = $headingRng.Find.Execute( `
    "1.3.2. Sustav za upravljanje klijentima (CRM)", $false, $false, $false, `
    $false, $false, $true, 1, $false, `
    "1.3.2. Sustav za upravljanje klijentima ", 2)

if (-not $found) {
    Write-Output "WARNING: CRM heading text was not found/replaced."
} else {
    Write-Output "Replaced CRM heading text."
}

$bookmarkInsertPos = $headingRng.End

# ---------------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark: delete it from its old spot (right after
#    the word "t" in "obuhvaćajući", near the top of the document) and add a
#    fresh zero-length bookmark right after the heading text we just edited.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
    Write-Output "Removed old _GoBack bookmark."
}
$bmRange = $d.Range($bookmarkInsertPos, $bookmarkInsertPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
Write-Output "Re-added _GoBack bookmark after the CRM heading."

# ---------------------------------------------------------------------------
# 3. Three empty "Normal (Web)" paragraphs only carry <w:numPr><w:numId val=0/>
#    and need an explicit <w:ilvl val="0"/> added ahead of the numId. Setting
#    ListLevelNumber stamps ilvl (but also forces a real numId), so immediately
#    strip the numbering back off with RemoveNumbers - that leaves ilvl=0 /
#    numId=0 behind, matching the target markup, while every other paragraph
#    property is left untouched.
# ---------------------------------------------------------------------------
$targetParaIndexes = @(30, 36, 40)
foreach ($idx in $targetParaIndexes) {
    $para = $d.Paragraphs($idx)
    $para.Range.ListFormat.ListLevelNumber = 1
    $para.Range.ListFormat.RemoveNumbers()
}
Write-Output "Stamped explicit ilvl=0 on the three empty list-less paragraphs."

# ---------------------------------------------------------------------------
# 4. Style "Normal (Web)" (internal id "8") gains an explicit rFonts entry
#    (matching the document's rPrDefault: Times New Roman / SimSun).
# ---------------------------------------------------------------------------
$normalWebFont = $d.Styles("Normal (Web)").Font
$normalWebFont.NameAscii = "Times New Roman"
$normalWebFont.NameOther = "Times New Roman"
$normalWebFont.NameFarEast = "SimSun"
$normalWebFont.NameBi = "Times New Roman"
Write-Output "Set explicit fonts on the 'Normal (Web)' style."

# ---------------------------------------------------------------------------
# 5. Style "Heading 2 Char" (internal id "11") becomes a quick style.
# ---------------------------------------------------------------------------
$d.Styles("Heading 2 Char").QuickStyle = $true
Write-Output "Marked 'Heading 2 Char' style as a quick style."
